$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = -3
    5  = 4
    7  = 2
    9  = 3
    10 = 5
    12 = -4
    13 = -6
    14 = 6
    15 = -4
    16 = 3
    17 = -2
    18 = 2
    19 = -7
    21 = -3
    22 = -1
    23 = 2
    24 = 4
    25 = 3
    26 = -1
    27 = 4
    28 = 1
    29 = 2
    31 = -2
    32 = 4
    33 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
